$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1729.0952
$ws.Range("I19").Value = 2053.5715
$ws.Range("J19").Value = 1566.8572
$ws.Range("K19").Value = 2053.5715
$ws.Range("L19").Value = 1566.8572
$ws.Range("M19").Value = -1878.5715
$ws.Range("N19").Value = -1916.8572
$ws.Range("H20").Value = 799
$ws.Range("I20").Value = 799
$ws.Range("K20").Value = 799
$ws.Range("M20").Value = -569
$ws.Range("H32").Value = 3796.6667
$ws.Range("I32").Value = 4440
$ws.Range("J32").Value = 2992.5
$ws.Range("K32").Value = 4440
$ws.Range("L32").Value = 2992.5
$ws.Range("M32").Value = -4114
$ws.Range("N32").Value = -3644.5
$ws.Range("H35").Value = 799
$ws.Range("I35").Value = 799
$ws.Range("K35").Value = 799
$ws.Range("M35").Value = -420
$ws.Range("H38").Value = 924.63635
$ws.Range("J38").Value = 5000
$ws.Range("L38").Value = 15000
$ws.Range("N38").Value = -15744
$ws.Range("H39").Value = 293.6154
$ws.Range("I39").Value = 234
$ws.Range("K39").Value = 702
$ws.Range("M39").Value = -406
$ws.Range("H43").Value = 5529
$ws.Range("I43").Value = 2498.3333
$ws.Range("J43").Value = 10075
$ws.Range("K43").Value = 2498.3333
$ws.Range("L43").Value = 10075
$ws.Range("M43").Value = -2429.3333
$ws.Range("N43").Value = -10213
$ws.Range("H58").Value = 172.5
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()  # was -7800
$ws.Range("H62").Value = 41499.332
$ws.Range("I62").Value = 31748.75
$ws.Range("K62").Value = 31748.75
$ws.Range("M62").Value = -31124.75
$ws.Range("H65").Value = 41499.332
$ws.Range("I65").Value = 31748.75
$ws.Range("K65").Value = 158743.75
$ws.Range("M65").Value = -155623.75
$ws.Range("H113").Value = 7362.875
$ws.Range("I113").Value = 5780.8
$ws.Range("K113").Value = 5780.8
$ws.Range("M113").Value = -2526.8
$ws.Range("H116").Value = 14326
$ws.Range("I116").Value = 13748.429
$ws.Range("K116").Value = 13748.429
$ws.Range("M116").Value = -10306.429
$ws.Range("H132").Value = 2688.8293
$ws.Range("I132").Value = 1253.973
$ws.Range("K132").Value = 3761.919
$ws.Range("M132").Value = -1231.919
$ws.Range("H137").Value = 3327.8333
$ws.Range("I137").Value = 2514.6
$ws.Range("J137").Value = 4141.067
$ws.Range("K137").Value = 7543.799999999999
$ws.Range("L137").Value = 12423.201
$ws.Range("M137").Value = -4993.799999999999
$ws.Range("N137").Value = -17523.201
$ws.Range("H138").Value = 2384.152
$ws.Range("I138").Value = 1571.3478
$ws.Range("J138").Value = 3196.9565
$ws.Range("K138").Value = 4714.0434
$ws.Range("L138").Value = 9590.869499999999
$ws.Range("M138").Value = 425.9565999999995
$ws.Range("N138").Value = -19870.8695

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 8769.5
$ws.Range("I53").Value = 8769.5
$ws.Range("K53").Value = 8769.5
$ws.Range("M53").Value = -8087.5
$ws.Range("H61").Value = 1936.75
$ws.Range("I61").Value = 1856.2858
$ws.Range("K61").Value = 1856.2858
$ws.Range("M61").Value = -1644.2858
$ws.Range("H88").Value = 10272.909
$ws.Range("J88").Value = 11185.5
$ws.Range("L88").Value = 11185.5
$ws.Range("N88").Value = -11997.5
$ws.Range("H91").Value = 10272.909
$ws.Range("J91").Value = 11185.5
$ws.Range("L91").Value = 11185.5
$ws.Range("N91").Value = -13993.5
$ws.Range("H124").Value = 56434.332
$ws.Range("J124").Value = 56434.332
$ws.Range("L124").Value = 56434.332
$ws.Range("N124").Value = -66254.33199999999
$ws.Range("H125").Value = 100715
$ws.Range("J125").Value = 100715
$ws.Range("L125").Value = 100715
$ws.Range("N125").Value = -110555
$ws.Range("H132").Value = 8299.625
$ws.Range("I132").Value = 8085.5
$ws.Range("J132").Value = 8942
$ws.Range("K132").Value = 24256.5
$ws.Range("L132").Value = 26826
$ws.Range("M132").Value = -21726.5
$ws.Range("N132").Value = -31886
$ws.Range("H136").Value = 1936.75
$ws.Range("I136").Value = 1856.2858
$ws.Range("K136").Value = 5568.857400000001
$ws.Range("M136").Value = -3018.857400000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 2400
$ws.Range("I42").Value = 2400
$ws.Range("K42").Value = 2400
$ws.Range("M42").Value = -1807
$ws.Range("H44").Value = 30749.75
$ws.Range("I44").Value = 26500
$ws.Range("J44").Value = 34999.5
$ws.Range("K44").Value = 26500
$ws.Range("L44").Value = 34999.5
$ws.Range("M44").Value = -26058
$ws.Range("N44").Value = -35883.5
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").ClearContents()  # was 10000
$ws.Range("N55").Value = 0
$ws.Range("H58").Value = 798.9091
$ws.Range("I58").Value = 689.2
$ws.Range("K58").Value = 689.2
$ws.Range("M58").Value = -486.2
$ws.Range("H93").Value = 26223.75
$ws.Range("I93").Value = 21999.5
$ws.Range("K93").Value = 21999.5
$ws.Range("M93").Value = -20127.5
$ws.Range("H134").Value = 1038.8695
$ws.Range("I134").Value = 836.8421
$ws.Range("J134").Value = 1998.5
$ws.Range("K134").Value = 2510.5263
$ws.Range("L134").Value = 5995.5
$ws.Range("M134").Value = 24.47370000000001
$ws.Range("N134").Value = -11065.5
$ws.Range("H136").Value = 798.9091
$ws.Range("I136").Value = 689.2
$ws.Range("K136").Value = 2067.6
$ws.Range("M136").Value = 482.3999999999996

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4030318.8
$ws.Range("I4").Value = 2918024.8
$ws.Range("K4").Value = 8754074.399999999
$ws.Range("M4").Value = -8753962.399999999
$ws.Range("H5").Value = 1456.8518
$ws.Range("I5").Value = 599.36365
$ws.Range("J5").Value = 2046.375
$ws.Range("K5").Value = 1798.09095
$ws.Range("L5").Value = 6139.125
$ws.Range("M5").Value = -1686.09095
$ws.Range("N5").Value = -6363.125
$ws.Range("H52").Value = 167026.33
$ws.Range("J52").Value = 431.8
$ws.Range("L52").Value = 1295.4
$ws.Range("N52").Value = -1827.4
$ws.Range("H80").Value = 8829925
$ws.Range("I80").Value = 1994
$ws.Range("K80").Value = 5982
$ws.Range("M80").Value = -5046
$ws.Range("H83").Value = 8829925
$ws.Range("I83").Value = 1994
$ws.Range("K83").Value = 17946
$ws.Range("M83").Value = -13266
$ws.Range("H97").Value = 372.66666
$ws.Range("I97").Value = 417.8
$ws.Range("J97").Value = 340.42856
$ws.Range("K97").Value = 1253.4
$ws.Range("L97").Value = 1021.28568
$ws.Range("M97").Value = -757.4000000000001
$ws.Range("N97").Value = -2013.28568
$ws.Range("H107").Value = 1024.8235
$ws.Range("I107").Value = 596.3333
$ws.Range("J107").Value = 1116.6428
$ws.Range("K107").Value = 1788.9999
$ws.Range("L107").Value = 3349.9284
$ws.Range("M107").Value = 131.0001
$ws.Range("N107").Value = -7189.928400000001
$ws.Range("H131").Value = 17316212
$ws.Range("J131").Value = 62330.176
$ws.Range("L131").Value = 186990.528
$ws.Range("N131").Value = -197070.528
$ws.Range("H135").Value = 1456.8518
$ws.Range("I135").Value = 599.36365
$ws.Range("J135").Value = 2046.375
$ws.Range("K135").Value = 5394.27285
$ws.Range("L135").Value = 18417.375
$ws.Range("M135").Value = -2859.27285
$ws.Range("N135").Value = -23487.375

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 36433
$ws.Range("J52").Value = 36433
$ws.Range("L52").Value = 36433
$ws.Range("N52").Value = -36951
$ws.Range("H80").Value = 8911.375
$ws.Range("I80").Value = 9776.308000000001
$ws.Range("J80").Value = 5163.3335
$ws.Range("K80").Value = 9776.308000000001
$ws.Range("L80").Value = 5163.3335
$ws.Range("M80").Value = -8778.308000000001
$ws.Range("N80").Value = -7159.3335
$ws.Range("H83").Value = 8911.375
$ws.Range("I83").Value = 9776.308000000001
$ws.Range("J83").Value = 5163.3335
$ws.Range("K83").Value = 48881.54000000001
$ws.Range("L83").Value = 25816.6675
$ws.Range("M83").Value = -43889.54000000001
$ws.Range("N83").Value = -35800.6675
$ws.Range("H128").Value = 53399.5
$ws.Range("J128").Value = 53399.5
$ws.Range("L128").Value = 53399.5
$ws.Range("N128").Value = -63359.5
$ws.Range("H132").Value = 7555.778
$ws.Range("I132").Value = 7555.778
$ws.Range("K132").Value = 22667.334
$ws.Range("M132").Value = -20137.334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 49070.285
$ws.Range("I7").Value = 49070.285
$ws.Range("K7").Value = 49070.285
$ws.Range("M7").Value = -48958.285
$ws.Range("H20").Value = 20000000
$ws.Range("J20").Value = 20000000
$ws.Range("L20").Value = 20000000
$ws.Range("N20").Value = -20000452
$ws.Range("H126").Value = 49070.285
$ws.Range("I126").Value = 49070.285
$ws.Range("K126").Value = 147210.855
$ws.Range("M126").Value = -144740.855
$ws.Range("H132").Value = 4995
$ws.Range("I132").Value = 4995
$ws.Range("K132").Value = 14985
$ws.Range("M132").Value = -12455
$ws.Range("H136").Value = 3043.7693
$ws.Range("I136").Value = 2406.9
$ws.Range("J136").Value = 5166.6665
$ws.Range("K136").Value = 7220.700000000001
$ws.Range("L136").Value = 15499.9995
$ws.Range("M136").Value = -4670.700000000001
$ws.Range("N136").Value = -20599.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5918.3335
$ws.Range("I132").Value = 5921
$ws.Range("J132").Value = 5897
$ws.Range("K132").Value = 17763
$ws.Range("L132").Value = 17691
$ws.Range("M132").Value = -15233
$ws.Range("N132").Value = -22751
$ws.Range("H136").Value = 6251.8823
$ws.Range("I136").Value = 4355.1333
$ws.Range("K136").Value = 13065.3999
$ws.Range("M136").Value = -10515.3999
